$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) column F needs updated values on both
# the "展览" and "全部类型" worksheets (rows 2-5).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 720
    $ws.Range("F3").Value = 4075
    $ws.Range("F4").Value = 114
    $ws.Range("F5").Value = 749
}
